# Refine metadata to be an additional tab
$wb = $excel.ActiveWorkbook

# --- Update time_taken (column F) timestamps on the "data" sheet ---
$dataSheet = $wb.Worksheets.Item("data")

$newTimes = @(
    "2021-10-05 14:33:32.096798",
    "2021-10-05 14:33:32.096806",
    "2021-10-05 14:33:32.096810",
    "2021-10-05 14:33:32.096812",
    "2021-10-05 14:33:32.096815",
    "2021-10-05 14:33:32.096818",
    "2021-10-05 14:33:32.096820",
    "2021-10-05 14:33:32.096823",
    "2021-10-05 14:33:32.096826",
    "2021-10-05 14:33:32.096829",
    "2021-10-05 14:33:32.096833",
    "2021-10-05 14:33:32.096836",
    "2021-10-05 14:33:32.096838",
    "2021-10-05 14:33:32.096841",
    "2021-10-05 14:33:32.096844",
    "2021-10-05 14:33:32.096846",
    "2021-10-05 14:33:32.096849",
    "2021-10-05 14:33:32.096852",
    "2021-10-05 14:33:32.096854",
    "2021-10-05 14:33:32.096857",
    "2021-10-05 14:33:32.096859",
    "2021-10-05 14:33:32.096862",
    "2021-10-05 14:33:32.096864",
    "2021-10-05 14:33:32.096867",
    "2021-10-05 14:33:32.096870",
    "2021-10-05 14:33:32.096872",
    "2021-10-05 14:33:32.096875",
    "2021-10-05 14:33:32.096878",
    "2021-10-05 14:33:32.096880",
    "2021-10-05 14:33:32.096883",
    "2021-10-05 14:33:32.096885",
    "2021-10-05 14:33:32.096888",
    "2021-10-05 14:33:32.096890",
    "2021-10-05 14:33:32.096893",
    "2021-10-05 14:33:32.096896",
    "2021-10-05 14:33:32.096898",
    "2021-10-05 14:33:32.096901"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- Add the new "metadata" sheet as an additional tab, placed after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Re-use the "data" sheet's existing header style (bold, centered, bordered)
# by copying its format onto the new sheet's header + index cells, so we
# don't fork a pile of new one-off styles.
$styleSource = $dataSheet.Cells.Item(1, 2)
$headerRng = $metaSheet.Range($metaSheet.Cells.Item(1, 2), $metaSheet.Cells.Item(1, 7))
$indexCell = $metaSheet.Cells.Item(2, 1)
$styleSource.Copy($headerRng)
$styleSource.Copy($indexCell)

# Header row
$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

# Data row
$metaSheet.Cells.Item(2, 1).Value = 0
$metaSheet.Cells.Item(2, 2).Value = "Cone-rod Dystrophy"
$metaSheet.Cells.Item(2, 3).Value = 3147
# data_version must stay textual ("0.31"), not be coerced to a number
$metaSheet.Cells.Item(2, 4).NumberFormat = "@"
$metaSheet.Cells.Item(2, 4).Value = "0.31"
$metaSheet.Cells.Item(2, 4).Style = "Normal"
$metaSheet.Cells.Item(2, 5).Value = "2021-09-18T08:16:33.459873Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:33:32.093205"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/3147/?format=json"
